# Auto-generated Excel COM-interop script to apply scheduled-runner market price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit worksheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 453.5
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1568
$ws.Range("H28").Value = 57073.39
$ws.Range("I28").Value = 57073.39
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 57073.39
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -56588.39
$ws.Range("N28").ClearContents()
$ws.Range("H40").Value = 6842.9287
$ws.Range("J40").Value = 12125
$ws.Range("L40").Value = 12125
$ws.Range("N40").Value = -12475
$ws.Range("H64").Value = 7235.294
$ws.Range("J64").Value = 7235.294
$ws.Range("L64").Value = 7235.294
$ws.Range("N64").Value = -7731.294
$ws.Range("H67").Value = 7235.294
$ws.Range("J67").Value = 7235.294
$ws.Range("L67").Value = 7235.294
$ws.Range("N67").Value = -8951.294
$ws.Range("H74").Value = 12318.0625
$ws.Range("I74").Value = 11006.923
$ws.Range("K74").Value = 11006.923
$ws.Range("M74").Value = -10070.923
$ws.Range("H77").Value = 12318.0625
$ws.Range("I77").Value = 11006.923
$ws.Range("K77").Value = 55034.61500000001
$ws.Range("M77").Value = -50354.61500000001
$ws.Range("H86").Value = 2634787.5
$ws.Range("I86").Value = 2075.111
$ws.Range("J86").Value = 4788825
$ws.Range("K86").Value = 2075.111
$ws.Range("L86").Value = 4788825
$ws.Range("M86").Value = -952.1109999999999
$ws.Range("N86").Value = -4791071
$ws.Range("H88").Value = 2020.7
$ws.Range("J88").Value = 2226.125
$ws.Range("L88").Value = 2226.125
$ws.Range("N88").Value = -3038.125
$ws.Range("H89").Value = 2634787.5
$ws.Range("I89").Value = 2075.111
$ws.Range("J89").Value = 4788825
$ws.Range("K89").Value = 10375.555
$ws.Range("L89").Value = 23944125
$ws.Range("M89").Value = -4759.555
$ws.Range("N89").Value = -23955357
$ws.Range("H91").Value = 2020.7
$ws.Range("J91").Value = 2226.125
$ws.Range("L91").Value = 2226.125
$ws.Range("N91").Value = -5034.125
$ws.Range("H106").Value = 3416
$ws.Range("I106").Value = 3284.4167
$ws.Range("K106").Value = 3284.4167
$ws.Range("M106").Value = -2653.4167
$ws.Range("H107").Value = 44508.566
$ws.Range("I107").Value = 48677.145
$ws.Range("K107").Value = 48677.145
$ws.Range("M107").Value = -46757.145
$ws.Range("H132").Value = 10578.435
$ws.Range("I132").Value = 1892.878
$ws.Range("J132").Value = 81800
$ws.Range("K132").Value = 5678.634
$ws.Range("L132").Value = 245400
$ws.Range("M132").Value = -3148.634
$ws.Range("N132").Value = -250460

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3542.481
$ws.Range("I32").Value = 3296.6624
$ws.Range("K32").Value = 3296.6624
$ws.Range("M32").Value = -3009.6624
$ws.Range("H61").Value = 3344.6072
$ws.Range("I61").Value = 3106.08
$ws.Range("J61").Value = 5332.3335
$ws.Range("K61").Value = 3106.08
$ws.Range("L61").Value = 5332.3335
$ws.Range("M61").Value = -2894.08
$ws.Range("N61").Value = -5756.3335
$ws.Range("H63").Value = 6644.4
$ws.Range("I63").Value = 3288.8
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 3288.8
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -2602.8
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 6644.4
$ws.Range("I66").Value = 3288.8
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 16444
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -13012
$ws.Range("N66").Value = -56864
$ws.Range("H74").Value = 1904.091
$ws.Range("I74").Value = 1741.375
$ws.Range("J74").Value = 2338
$ws.Range("K74").Value = 1741.375
$ws.Range("L74").Value = 2338
$ws.Range("M74").Value = -867.375
$ws.Range("N74").Value = -4086
$ws.Range("H77").Value = 1904.091
$ws.Range("I77").Value = 1741.375
$ws.Range("J77").Value = 2338
$ws.Range("K77").Value = 8706.875
$ws.Range("L77").Value = 11690
$ws.Range("M77").Value = -4338.875
$ws.Range("N77").Value = -20426
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H136").Value = 3344.6072
$ws.Range("I136").Value = 3106.08
$ws.Range("J136").Value = 5332.3335
$ws.Range("K136").Value = 9318.24
$ws.Range("L136").Value = 15997.0005
$ws.Range("M136").Value = -6768.24
$ws.Range("N136").Value = -21097.0005

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 39128.57
$ws.Range("I82").Value = 6780
$ws.Range("K82").Value = 6780
$ws.Range("M82").Value = -6397
$ws.Range("H85").Value = 39128.57
$ws.Range("I85").Value = 6780
$ws.Range("K85").Value = 6780
$ws.Range("M85").Value = -5454
$ws.Range("H86").Value = 898450.6
$ws.Range("I86").Value = 1549104.1
$ws.Range("K86").Value = 1549104.1
$ws.Range("M86").Value = -1547981.1
$ws.Range("H89").Value = 898450.6
$ws.Range("I89").Value = 1549104.1
$ws.Range("K89").Value = 7745520.5
$ws.Range("M89").Value = -7739904.5
$ws.Range("H105").Value = 2950.25
$ws.Range("I105").Value = 2950.25
$ws.Range("K105").Value = 2950.25
$ws.Range("M105").Value = -1203.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3495.125
$ws.Range("I62").Value = 2509.1667
$ws.Range("K62").Value = 2509.1667
$ws.Range("M62").Value = -1885.1667
$ws.Range("H65").Value = 3495.125
$ws.Range("I65").Value = 2509.1667
$ws.Range("K65").Value = 12545.8335
$ws.Range("M65").Value = -9425.833500000001
$ws.Range("H107").Value = 479.46155
$ws.Range("I107").Value = 409.2857
$ws.Range("K107").Value = 409.2857
$ws.Range("M107").Value = 1510.7143
$ws.Range("H134").Value = 296545.1
$ws.Range("I134").Value = 2501
$ws.Range("K134").Value = 7503
$ws.Range("M134").Value = -4968

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1544018.8
$ws.Range("I80").Value = 1254907.4
$ws.Range("K80").Value = 1254907.4
$ws.Range("M80").Value = -1253909.4
$ws.Range("H83").Value = 1544018.8
$ws.Range("I83").Value = 1254907.4
$ws.Range("K83").Value = 6274537
$ws.Range("M83").Value = -6269545
$ws.Range("H110").Value = 40117
$ws.Range("J110").Value = 40117
$ws.Range("L110").Value = 40117
$ws.Range("N110").Value = -48297
$ws.Range("H113").Value = 772345.0600000001
$ws.Range("I113").Value = 2000940.4
$ws.Range("J113").Value = 4473
$ws.Range("K113").Value = 2000940.4
$ws.Range("L113").Value = 4473
$ws.Range("M113").Value = -1998770.4
$ws.Range("N113").Value = -8813

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8248.375
$ws.Range("I7").Value = 8458.4
$ws.Range("K7").Value = 8458.4
$ws.Range("M7").Value = -8346.4
$ws.Range("H46").Value = 2631.7273
$ws.Range("I46").Value = 2583
$ws.Range("K46").Value = 2583
$ws.Range("M46").Value = -2395
$ws.Range("H68").Value = 2722
$ws.Range("I68").Value = 1444
$ws.Range("K68").Value = 1444
$ws.Range("M68").Value = -695
$ws.Range("H71").Value = 2722
$ws.Range("I71").Value = 1444
$ws.Range("K71").Value = 7220
$ws.Range("M71").Value = -3476
$ws.Range("H122").Value = 6779.4
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 8248.375
$ws.Range("I126").Value = 8458.4
$ws.Range("K126").Value = 25375.2
$ws.Range("M126").Value = -22905.2
$ws.Range("H136").Value = 232942.5
$ws.Range("I136").Value = 374150.4
$ws.Range("J136").Value = 8671.117
$ws.Range("K136").Value = 1122451.2
$ws.Range("L136").Value = 26013.351
$ws.Range("M136").Value = -1119901.2
$ws.Range("N136").Value = -31113.351

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9594765
$ws.Range("I136").Value = 12838487
$ws.Range("J136").Value = 224015.56
$ws.Range("K136").Value = 38515461
$ws.Range("L136").Value = 672046.6799999999
$ws.Range("M136").Value = -38512911
$ws.Range("N136").Value = -677146.6799999999

Write-Output "Applied scheduled Sheets update: $( 208 ) value writes, $( 3 ) cell clears across 7 worksheets."